# The committed change refreshes the "quadratic-svm-score" sheet with a new
# run's numeric scores in column B (the "1-c__Elusimicrobia" score column);
# everything else (headers, row labels, column C) is unchanged.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1643.9452004067166
$ws.Range("B3").Value = 1535.0922864223048
$ws.Range("B4").Value = 1723.472007440319
